# Suivi sheet: a new price-history snapshot column was inserted right
# before the "nom" column (previously column DW), shifting "nom" and
# "url_produit" one column to the right (DW->DX, DX->DY).
#
# New column DW:
#   - Row 1 (header): new timestamp "2026-02-02 13:54:36"
#   - Rows 2-206: the new price snapshot. For rows that still had a
#     price in the previous snapshot column (now DV), the price is
#     carried forward unchanged; rows whose price history had already
#     gone blank stay blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at DW (column 127), shifting "nom"/"url_produit"
# (and their column-1 style) one column to the right, just like Excel's
# Insert Sheet Columns command.
$ws.Columns.Item(127).Insert()

# New header timestamp cell.
$ws.Cells.Item(1, 127).Value = "2026-02-02 13:54:36"

# Carry forward the latest known price (from column DV, column 126) into
# the freshly inserted column DW for every row that still had a value.
for ($row = 2; $row -le 80; $row++) {
    $price = $ws.Cells.Item($row, 126).Value()
    $ws.Cells.Item($row, 127).Value = $price
}
